# Update the "removeCouple.php" row (row 16):
#   - Parameters to send (B16): "token" -> "token, token2"
#   - What is for (E16): append clarification about Token2
# Update the "updateCouple.php" row (row 29):
#   - Data returned (C29): append the new "already_paired" response case

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "token, token2"
$ws.Range("E16").Value = "remove current couple from specified token. Token2 is the other side person."

$ws.Range("C29").Value = "[success:0] [success:1, img] [`"success`":`"already_paired`"]"

# Leave the view scrolled to where the edits happened, matching the author's
# last on-screen selection when the file was saved.
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 25
